$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 41.81968580500756
$ws.Range("C2").Value = 19.72484728363376
$ws.Range("D2").Value = 11.12811060419194
$ws.Range("E2").Value = 9.042288778937897
$ws.Range("G2").Value = 3.929090748422566
$ws.Range("J2").Value = 6.676064186320634
$ws.Range("L2").Value = 14.80184051889697
$ws.Range("M2").Value = 29.71305897150825
$ws.Range("N2").Value = 29.5381333998286

$ws.Range("B3").Value = 41.62657699828047
$ws.Range("C3").Value = 19.46864753596805
$ws.Range("D3").Value = 11.15246419390172
$ws.Range("E3").Value = 9.028133989954034
$ws.Range("G3").Value = 3.935420095840056
$ws.Range("J3").Value = 6.684976690597018
$ws.Range("L3").Value = 14.8296264514482
$ws.Range("M3").Value = 29.70262071842842
$ws.Range("N3").Value = 29.46644932299319

$ws.Range("B4").Value = 41.51913942877939
$ws.Range("C4").Value = 19.31638728111385
$ws.Range("D4").Value = 11.16900268920494
$ws.Range("E4").Value = 9.019294785846288
$ws.Range("G4").Value = 3.939498984438616
$ws.Range("J4").Value = 6.691722240149192
$ws.Range("L4").Value = 14.84909188408768
$ws.Range("M4").Value = 29.70279658987139
$ws.Range("N4").Value = 29.4234364964641

$ws.Range("B5").Value = 41.47817857120293
$ws.Range("C5").Value = 19.25567685958344
$ws.Range("D5").Value = 11.1761407195411
$ws.Range("E5").Value = 9.015655014424995
$ws.Range("G5").Value = 3.941209847264974
$ws.Range("J5").Value = 6.694790398198344
$ws.Range("L5").Value = 14.85762797070824
$ws.Range("M5").Value = 29.70452133831072
$ws.Range("N5").Value = 29.40616152897514

$ws.Range("B6").Value = 41.47154797351559
$ws.Range("C6").Value = 19.24567859592596
$ws.Range("D6").Value = 11.17735004479958
$ws.Range("E6").Value = 9.015048340653737
$ws.Range("G6").Value = 3.941496882092795
$ws.Range("J6").Value = 6.695319121287261
$ws.Range("L6").Value = 14.85908181930355
$ws.Range("M6").Value = 29.70490747541375
$ws.Range("N6").Value = 29.4033083543685

$ws.Range("B7").Value = 41.51857556935373
$ws.Range("C7").Value = 19.31556301922101
$ws.Range("D7").Value = 11.1690973420526
$ws.Range("E7").Value = 9.0192458519965
$ws.Range("G7").Value = 3.939521860298498
$ws.Range("J7").Value = 6.691762326788014
$ws.Range("L7").Value = 14.84920456135603
$ws.Range("M7").Value = 29.70281316143327
$ws.Range("N7").Value = 29.42320249384055

$ws.Range("B8").Value = 41.75080560211
$ws.Range("C8").Value = 19.63550356488874
$ws.Range("D8").Value = 11.13617858959217
$ws.Range("E8").Value = 9.037438654315503
$ws.Range("G8").Value = 3.931233269706005
$ws.Range("J8").Value = 6.678872447628338
$ws.Range("L8").Value = 14.81092145083797
$ws.Range("M8").Value = 29.70809172286428
$ws.Range("N8").Value = 29.51320967721567

$ws.Range("B9").Value = 42.29340582060071
$ws.Range("C9").Value = 20.2998420610836
$ws.Range("D9").Value = 11.08421200836722
$ws.Range("E9").Value = 9.071962369047377
$ws.Range("G9").Value = 3.916496902959245
$ws.Range("J9").Value = 6.663735539809862
$ws.Range("L9").Value = 14.75497055215338
$ws.Range("M9").Value = 29.77078369121423
$ws.Range("N9").Value = 29.69760336415244

$ws.Range("B10").Value = 42.74345007531227
$ws.Range("C10").Value = 20.80613798154306
$ws.Range("D10").Value = 11.05371787231112
$ws.Range("E10").Value = 9.096657744773493
$ws.Range("G10").Value = 3.906579777276611
$ws.Range("J10").Value = 6.658845288941557
$ws.Range("L10").Value = 14.72557740914103
$ws.Range("M10").Value = 29.84883063461644
$ws.Range("N10").Value = 29.8378116559472

$ws.Range("B11").Value = 42.95889677235646
$ws.Range("C11").Value = 21.03948064650939
$ws.Range("D11").Value = 11.04151684371811
$ws.Range("E11").Value = 9.107754301821046
$ws.Range("G11").Value = 3.902262415550502
$ws.Range("J11").Value = 6.657982874647367
$ws.Range("L11").Value = 14.71476069161375
$ws.Range("M11").Value = 29.89128035756569
$ws.Range("N11").Value = 29.90260469906255

$ws.Range("B12").Value = 43.04197907435984
$ws.Range("C12").Value = 21.12819802861126
$ws.Range("D12").Value = 11.03713711638122
$ws.Range("E12").Value = 9.11193706523485
$ws.Range("G12").Value = 3.900655176522969
$ws.Range("J12").Value = 6.657852808765279
$ws.Range("L12").Value = 14.7110328287712
$ws.Range("M12").Value = 29.90835184350541
$ws.Range("N12").Value = 29.92728337942351

$ws.Range("B13").Value = 43.02401995359516
$ws.Range("C13").Value = 21.10907659548122
$ws.Range("D13").Value = 11.0380696669338
$ws.Range("E13").Value = 9.11103708081486
$ws.Range("G13").Value = 3.901000098248485
$ws.Range("J13").Value = 6.657872071558699
$ws.Range("L13").Value = 14.71181930178004
$ws.Range("M13").Value = 29.90463090845673
$ws.Range("N13").Value = 29.92196206746702

$ws.Range("B14").Value = 42.96570225791588
$ws.Range("C14").Value = 21.04677280505723
$ws.Range("D14").Value = 11.04115169869889
$ws.Range("E14").Value = 9.108098803381827
$ws.Range("G14").Value = 3.902129634235084
$ws.Range("J14").Value = 6.657968232295694
$ws.Range("L14").Value = 14.71444661410281
$ws.Range("M14").Value = 29.89266487903819
$ws.Range("N14").Value = 29.90463217877762

$ws.Range("B15").Value = 42.93017462103423
$ws.Range("C15").Value = 21.0086538951844
$ws.Range("D15").Value = 11.04307086492849
$ws.Range("E15").Value = 9.106296530361593
$ws.Range("G15").Value = 3.902825101533704
$ws.Range("J15").Value = 6.658052742685445
$ws.Range("L15").Value = 14.71610389359946
$ws.Range("M15").Value = 29.88546506167057
$ws.Range("N15").Value = 29.89403566196685

$ws.Range("B16").Value = 42.72958212366332
$ws.Range("C16").Value = 20.79094256203846
$ws.Range("D16").Value = 11.0545488837527
$ws.Range("E16").Value = 9.095929835390454
$ws.Range("G16").Value = 3.906865809800785
$ws.Range("J16").Value = 6.658929115931294
$ws.Range("L16").Value = 14.72633578372379
$ws.Range("M16").Value = 29.84619605428112
$ws.Range("N16").Value = 29.83359742390962

$ws.Range("B17").Value = 42.60924016061856
$ws.Range("C17").Value = 20.65810191888265
$ws.Range("D17").Value = 11.06201842474612
$ws.Range("E17").Value = 9.089535349367756
$ws.Range("G17").Value = 3.909394165845525
$ws.Range("J17").Value = 6.659816101933124
$ws.Range("L17").Value = 14.7332675929169
$ws.Range("M17").Value = 29.82388357565884
$ws.Range("N17").Value = 29.79677788910587

$ws.Range("B18").Value = 42.54103443707526
$ws.Range("C18").Value = 20.58198505835169
$ws.Range("D18").Value = 11.06647196549892
$ws.Range("E18").Value = 9.085844437849723
$ws.Range("G18").Value = 3.910866683086393
$ws.Range("J18").Value = 6.660454454847735
$ws.Range("L18").Value = 14.73749496640857
$ws.Range("M18").Value = 29.81170398314454
$ws.Range("N18").Value = 29.77569564340567

$ws.Range("B19").Value = 42.51811622211834
$ws.Range("C19").Value = 20.55626536745984
$ws.Range("D19").Value = 11.06800685902617
$ws.Range("E19").Value = 9.084592515203857
$ws.Range("G19").Value = 3.911368398038544
$ws.Range("J19").Value = 6.660692585177451
$ws.Range("L19").Value = 14.73896754036012
$ws.Range("M19").Value = 29.80769254821672
$ws.Range("N19").Value = 29.76857401160591

$ws.Range("B20").Value = 42.6219463617452
$ws.Range("C20").Value = 20.67221364245174
$ws.Range("D20").Value = 11.06120700179645
$ws.Range("E20").Value = 9.090217392227681
$ws.Range("G20").Value = 3.909123128527356
$ws.Range("J20").Value = 6.659708409617536
$ws.Range("L20").Value = 14.73250480504861
$ws.Range("M20").Value = 29.82619110747584
$ws.Range("N20").Value = 29.80068753788554

$ws.Range("B21").Value = 42.98279130798356
$ws.Range("C21").Value = 21.0650639178131
$ws.Range("D21").Value = 11.04023990107639
$ws.Range("E21").Value = 9.108962365367185
$ws.Range("G21").Value = 3.901797113605542
$ws.Range("J21").Value = 6.657934649600832
$ws.Range("L21").Value = 14.71366490937465
$ws.Range("M21").Value = 29.89615256319768
$ws.Range("N21").Value = 29.9097185280522

$ws.Range("B22").Value = 43.22732939335394
$ws.Range("C22").Value = 21.32385074217938
$ws.Range("D22").Value = 11.02793884283759
$ws.Range("E22").Value = 9.12110125617863
$ws.Range("G22").Value = 3.897170203702787
$ws.Range("J22").Value = 6.65792105886476
$ws.Range("L22").Value = 14.70349843281519
$ws.Range("M22").Value = 29.94768430309261
$ws.Range("N22").Value = 29.98180858631494

$ws.Range("B23").Value = 43.09603360696578
$ws.Range("C23").Value = 21.18557102746519
$ws.Range("D23").Value = 11.03437576643337
$ws.Range("E23").Value = 9.114632582965598
$ws.Range("G23").Value = 3.899625015800409
$ws.Range("J23").Value = 6.657823293681647
$ws.Range("L23").Value = 14.7087277920689
$ws.Range("M23").Value = 29.91965035523804
$ws.Range("N23").Value = 29.94325750908726

$ws.Range("B24").Value = 42.61619883094806
$ws.Range("C24").Value = 20.66583293192809
$ws.Range("D24").Value = 11.06157335025286
$ws.Range("E24").Value = 9.089909086278089
$ws.Range("G24").Value = 3.909245605512405
$ws.Range("J24").Value = 6.659756697384618
$ws.Range("L24").Value = 14.73284890716178
$ws.Range("M24").Value = 29.82514585335498
$ws.Range("N24").Value = 29.79891971825119

$ws.Range("B25").Value = 42.137450128583
$ws.Range("C25").Value = 20.116606233903
$ws.Range("D25").Value = 11.09692152646299
$ws.Range("E25").Value = 9.062742724471045
$ws.Range("G25").Value = 3.920322612248976
$ws.Range("J25").Value = 6.666739577155578
$ws.Range("L25").Value = 14.76805356871254
$ws.Range("M25").Value = 29.74820862432848
$ws.Range("N25").Value = 29.6468823492705
